# Plan de sprint aleph.xlsx -- "Correccion de stories y se colocó tiempos
# reales de los stories terminados"
#
# Sheet "Plan Sprint v1.1" (index 2): just a selection/cursor change.
# Sheet "Plan Sprint v1.2" (index 3): fill in the real hours / assignment
# data for a few stories that were finished (rows 5, 6 and 8), record the
# "Real" hours for rows 3 and 4, shrink row 4's height, and move the
# active cell / cursor.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Plan Sprint v1.1" -- only the selected range changed.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("H3:N8").Select()

# ---------------------------------------------------------------------
# "Plan Sprint v1.2" -- story corrections + real hours.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()

# Row 3 (A1 / Login): real hours puntuation + "Terminado" = ok
$ws3.Range("S3").Value = 12
$ws3.Range("U3").Value = "ok"

# Row 4 (A2 / Alta de Usuarios): real hours puntuation + "Terminado" = ok
# also shrink the custom row height from 53 to 46
$ws3.Range("S4").Value = 8
$ws3.Range("U4").Value = "ok"
$ws3.Rows.Item(4).RowHeight = 46

# Row 5 (A3 / Estudios Asignados, finished by Miguel -> now N1 / Estudios
# Asignados, A, 5/0, ok)
$ws3.Range("H5").Value = "N1"
$ws3.Range("I5").Value = "Estudios Asignados"
$ws3.Range("J5").Value = "A"
$ws3.Range("K5").Value = 5
$ws3.Range("L5").Value = 0
$ws3.Range("M5").Value = "ok"

# Row 6 (N2 / Estudios En Particular, A, 8/0, David)
$ws3.Range("H6").Value = "N2"
$ws3.Range("I6").Value = "Estudios En Particular"
$ws3.Range("J6").Value = "A"
$ws3.Range("K6").Value = 8
$ws3.Range("L6").Value = 0
$ws3.Range("M6").Value = "David"

# Row 8 (L1 / Estudios En Particular, A, 2/0, Rodrigo)
$ws3.Range("H8").Value = "L1"
$ws3.Range("I8").Value = "Estudios En Particular"
$ws3.Range("J8").Value = "A"
$ws3.Range("K8").Value = 2
$ws3.Range("L8").Value = 0
$ws3.Range("M8").Value = "Rodrigo"

# Move the cursor / active cell to where the editor left off.
$ws3.Range("P4").Select()
$excel.ActiveWindow.ScrollColumn = 13
$excel.ActiveWindow.ScrollRow = 1
